$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scaling")

# Update the H column (C1) values in rows 3-17 to sequential numbers 2..16
for ($row = 3; $row -le 17; $row++) {
    $ws.Cells.Item($row, 8).Value = $row - 1
}

# Update the active selection to J12, matching the edited sheet view
$ws.Range("J12").Select()
